# Update "Saroum Cement" yearly balance-sheet workbook:
#  - shift the 5 reporting periods / publish dates one column to the left
#    and append the new (6th) period on the right
#  - shift every financial figure one column to the left the same way,
#    appending the freshly reported period's figures on the right
#  - a handful of cells toggle between a literal "0" and the "-" placeholder
#    used by the source for "not reported" figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 8: "دوره مالی" (financial period) headers
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# ---------------------------------------------------------------------------
# Row 9: "تاریخ انتشار" (publish date) headers
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = "1399-01-24 (7)"
$ws.Range("E9").Value = "1400-02-04 (7)"
$ws.Range("F9").Value = "1401-01-31 (8)"
$ws.Range("G9").Value = "1402-01-30 (9)"
$ws.Range("H9").Value = "1402-01-30 (2)"

# ---------------------------------------------------------------------------
# Helper-style explicit per-row values (D:H) for the data rows
# ---------------------------------------------------------------------------
$rows = @{
    12 = @(85481, 443530, 117259, 253515, 435412)
    13 = @(0, 84055, 368800, 0, 0)
    14 = @(390177, 337359, 464357, 385305, 560921)
    15 = @(676551, 810986, 1142222, 2400230, 3175528)
    16 = @(20672, 55240, 583408, 270598, 354677)
    17 = @(0, 0, 0, 0, 0)
    18 = @(1172881, 1731170, 2676046, 3309648, 4526538)
    19 = @(4928, 2750, 1844, 13393, 46815)
    20 = @(11515, 9116, 946405, 3060535, 2891371)
    21 = @(0, 0, 0, 0, 0)
    22 = @(783918, 637525, 579796, 939707, 2186586)
    23 = @(46989, 46994, 46995, 46995, 46995)
    25 = @(0, 0, 0, 0, 0)
    26 = @(847350, 696385, 1575040, 4060630, 5171767)
    27 = @(2020231, 2427555, 4251086, 7370278, 9698305)
    29 = @(159904, 127572, 229156, 667745, 678699)
    31 = @(70858, 37451, 163348, 544501, 713511)
    32 = @(84040, 132039, 202252, 394434, 483812)
    33 = @(192023, 22463, 38693, 57542, 96131)
    34 = @(121235, 121235, 121235, 0, 358975)
    35 = @(0, 13350, 0, 0, 0)
    36 = @(0, 0, 0, 0, 0)
    37 = @(628060, 454110, 754684, 1664222, 2331128)
    38 = @(0, 0, 0, 0, 0)
    40 = @(0, 0, 0, 121235, 121235)
    41 = @(49605, 64202, 87343, 128883, 213522)
    42 = @(49605, 64202, 87343, 250118, 334757)
    43 = @(677665, 518312, 842027, 1914340, 2665885)
    45 = @(700000, 700000, 700000, 700000, 1000000)
    46 = @(0, 0, 0, 0, 0)
    47 = @(0, 0, 0, 0, 0)
    48 = @(0, 0, -57509, -80960, -89021)
    49 = @(0, 0, 0, 0, 32129)
    50 = @(70000, 70000, 70000, 70000, 100000)
    51 = @(0, 0, 0, 0, 0)
    53 = @(0, 0, 0, 0, 0)
    55 = @(0, 0, 0, 0, 0)
    56 = @(572566, 1139243, 2696568, 4766898, 5989312)
    57 = @(1342566, 1909243, 3409059, 5455938, 7032420)
    58 = @(2020231, 2427555, 4251086, 7370278, 9698305)
}

$cols = @("D", "E", "F", "G", "H")

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------------
# Cells that flip between the literal 0 and the "-" (not reported) marker
# ---------------------------------------------------------------------------
$ws.Range("D39").Value = "-"
$ws.Range("D49").Value = 0
$ws.Range("D52").Value = "-"
$ws.Range("D54").Value = "-"
